$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# --- Paragraph 1 (University of Washington) ---------------------------------
# Text content is unchanged, but re-typing the paragraph (as happened in the
# source edit) collapses the old run fragments / spell-check proofErr marks
# into a single clean run.
$p1Text = "University of Washington: wpneuro (C: 7/26, R: 8/11), Bear08150902 (C: 7/14, R: 7/26), JimmyMacElroy (C: 7/20, R: 8/12), neurofanatic213 (C: 7/24, R: 8/12), supersonic99 (C: 7/21, R: 8/15), gnohmChild (C: 7/14, R: 7/16), DCemt (C: 7/28, R:8/11), scineurd (C: 7/24, R: 7/28), scientistintraining (C: 7/25, R: 8/12), ac_505 (C: 8/24, R: 9/6), wrhen (C: 9/X, R: 9/22)"
$find.Execute($p1Text, $true, $false, $false, $false, $false, $true, 1, $false, $p1Text, 2) | Out-Null

# --- Paragraph 2 (University of California San Francisco) -------------------
$p2Text = "University of California San Francisco: dysyurt (C: 9/30, R: 12/3)"
$find.Execute($p2Text, $true, $false, $false, $false, $false, $true, 1, $false, $p2Text, 2) | Out-Null

# --- Paragraph 4 (Emory University) ------------------------------------------
# Fix the typo "R, 9/25" -> "R: 9/25" (comma -> colon) by editing only the
# single punctuation character, the way a user would in Word: select just
# that character and retype it. Locate it with Find so the edit is robust to
# any earlier shifts in the document.
$emoryRange = $d.Content
$emoryRange.Find.ClearFormatting()
$found = $emoryRange.Find.Execute("Emory University: neurozf5 (C: 8/3, R, 9/25)")
if ($found) {
    $fullText = $emoryRange.Text
    $commaOffset = $fullText.IndexOf(", 9/25")
    $commaStart = $emoryRange.Start + $commaOffset
    $commaRange = $d.Range($commaStart, $commaStart + 1)
    $commaRange.Text = ":"
    # Nudge formatting on just the new character so Word keeps it as its own
    # run instead of silently re-merging it with its neighbors.
    $commaRange.Bold = 1
    $commaRange.Bold = 0
}
